$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewPolicy")

# Row 3: new policy entry "newpol_002" / "P", mirroring the pattern of row 2
$ws.Range("A3").Value = "newpol_002"
$ws.Range("B3").Value = "P"

# D3 must become an (empty-string) text cell, matching D2's existing empty text value.
# A plain Value="" clears the cell instead of storing an empty string, so we force
# text-typing via the quote-prefix entry trick, then strip the formatting residue.
$ws.Range("D3").Formula = "'"
$ws.Range("D3").Style = "Normal"

# Move the active selection to B2
$ws.Range("B2").Select()
